$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.370.98"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").Value = "3.116.49"
$ws.Range("E3").Value = "  -2.74%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +1.18%  "

$ws.Range("D5").Value = "'561.53"
$ws.Range("E5").Value = "  -4.09%  "

$ws.Range("D6").Value = "'135.10"
$ws.Range("E6").Value = "  -9.36%  "

$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.61%  "

$ws.Range("D8").Value = "3.126.62"
$ws.Range("E8").Value = "  -2.21%  "

$ws.Range("D9").Value = "'0.506"
$ws.Range("E9").Value = "  -5.20%  "

$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  -8.55%  "

$ws.Range("D11").Value = "'5.94"
$ws.Range("E11").Value = "  -5.07%  "

$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  -3.81%  "

$ws.Range("D13").Value = "'0.0000220"
$ws.Range("E13").Value = "  -7.86%  "

$ws.Range("D14").Value = "'34.06"
$ws.Range("E14").Value = "  -10.86%  "

$ws.Range("D15").Value = "3.681.36"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "65.885.01"
$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("D17").Value = "3.188.66"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").Value = "'6.56"
$ws.Range("E19").Value = "  -4.33%  "

$ws.Range("D20").Value = "'482.30"
$ws.Range("E20").Value = "  -6.80%  "

$ws.Range("D21").Value = "'13.86"
$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").Value = "'0.688"
$ws.Range("E22").Value = "  -6.97%  "

$ws.Range("D23").Value = "'7.05"
$ws.Range("E23").Value = "  -9.08%  "

$ws.Range("D24").Value = "'78.93"
$ws.Range("E24").Value = "  -5.65%  "

$ws.Range("D25").Value = "'12.31"
$ws.Range("E25").Value = "  -6.69%  "

$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").Value = "'2.97"
$ws.Range("E27").Value = "  -5.20%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'1.94"
$ws.Range("E28").Value = "  -7.19%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'26.49"
$ws.Range("E29").Value = "  -6.71%  "

$ws.Range("D30").Value = "'7.19"
$ws.Range("E30").Value = "  -7.66%  "

$ws.Range("D31").Value = "'1.14"
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").Value = "'2.43"
$ws.Range("E32").Value = "  -6.98%  "

$ws.Range("D33").Value = "'1.01"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").Value = "'490.80"
$ws.Range("E34").Value = "  -7.86%  "

$ws.Range("D35").Value = "'53.45"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").Value = "'5.77"
$ws.Range("E36").Value = "  -8.21%  "

$ws.Range("D37").Value = "'5.02"
$ws.Range("E37").Value = "  -10.26%  "

$ws.Range("D38").Value = "'0.0397"
$ws.Range("E38").Value = "  -5.99%  "

$ws.Range("D39").Value = "'0.0781"
$ws.Range("E39").Value = "  -7.52%  "

$ws.Range("D40").Value = "'8.27"
$ws.Range("E40").Value = "  -8.64%  "

$ws.Range("D41").Value = "'0.114"
$ws.Range("E41").Value = "  -7.17%  "

$ws.Range("D42").Value = "2.776.79"
$ws.Range("E42").Value = "  -2.81%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.35"
$ws.Range("E44").Value = "  -13.50%  "

$ws.Range("D45").Value = "'0.238"
$ws.Range("E45").Value = "  -6.60%  "

$ws.Range("D46").Value = "'120.17"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.105"
$ws.Range("E47").Value = "  -6.54%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'23.45"
$ws.Range("E48").Value = "  -8.27%  "

$ws.Range("D49").Value = "'1.90"
$ws.Range("E49").Value = "  -10.03%  "

$ws.Range("D50").Value = "0.0$([char]0x2083)0488"
$ws.Range("E50").Value = "  -14.82%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.02"
$ws.Range("E51").Value = "  -14.87%  "
